$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 24 data
$ws.Cells.Item(24,1).Value2 = 45866
$ws.Cells.Item(24,1).NumberFormat = $ws.Cells.Item(23,1).NumberFormat
$ws.Cells.Item(24,2).Value2 = 6389.77001953125
$ws.Cells.Item(24,3).Value2 = 6401.06982421875
$ws.Cells.Item(24,4).Value2 = 6375.7900390625
$ws.Cells.Item(24,5).Value2 = 6397.68994140625
$ws.Cells.Item(24,6).Value2 = 4565620000
$ws.Cells.Item(24,7).Value2 = 0.0001768581088181
$ws.Cells.Item(24,8).Value2 = 1
$ws.Cells.Item(24,9).Value2 = 'neon bloom inc announces successful launch of aidriven backend software facebook parent meta to halt political ad in eu citing unworkable new regulation education department say it will release billion in remaining withheld grant money for school education department say it will release billion in remaining withheld grant money for school education department say it will release billion in remaining withheld grant money for school education department say it will release billion in remaining withheld grant money for school rd place at southbridge celebrates retail anchor with grand opening ribboncutting tempus wealth planning llc ha position in alphabet inc nasdaqgoog neo ivy capital management cut stock holding in johnson johnson nysejnj nicholas wealth llc trim position in alphabet inc nasdaqgoog aldi trader joes adding store in pennsylvania here where when cocoa price sharply higher on belownormal rain in west africa everything you need to know about relocating part everything you need to know about relocating part everything you need to know about relocating part everything you need to know about relocating part palantir join list of most valuable u company with stock more than doubling in coffee price slump a brazil coffee harvest accelerates chicago police blew past dnc overtime budget foreign investor in charlotte eagle lake jf project obtain conditional green card approval it will be a tidal wave of earning and data next week ap sport summarybrief at pm edt ap sport summarybrief at pm edt nearshoring a smart way to cut cost and scale your business in nearshoring a smart way to cut cost and scale your business in nearshoring a smart way to cut cost and scale your business in nearshoring a smart way to cut cost and scale your business in nearshoring a smart way to cut cost and scale your business in nearshoring a smart way to cut cost and scale your business in nearshoring a smart way to cut cost and scale your business in nearshoring a smart way to cut cost and scale your business in nearshoring a smart way to cut cost and scale your business in nearshoring a smart way to cut cost and scale your business in nearshoring a smart way to cut cost and scale your business in nearshoring a smart way to cut cost and scale your business in nearshoring a smart way to cut cost and scale your business in nearshoring a smart way to cut cost and scale your business in nearshoring a smart way to cut cost and scale your business in nearshoring a smart way to cut cost and scale your business in nearshoring a smart way to cut cost and scale your business in nearshoring a smart way to cut cost and scale your business in nearshoring a smart way to cut cost and scale your business in nearshoring a smart way to cut cost and scale your business in nearshoring a smart way to cut cost and scale your business in nearshoring a smart way to cut cost and scale your business in nearshoring a smart way to cut cost and scale your business in nearshoring a smart way to cut cost and scale your business in nearshoring a smart way to cut cost and scale your business in nearshoring a smart way to cut cost and scale your business in nearshoring a smart way to cut cost and scale your business in nearshoring a smart way to cut cost and scale your business in nearshoring a smart way to cut cost and scale your business in nearshoring a smart way to cut cost and scale your business in nearshoring a smart way to cut cost and scale your business in nearshoring a smart way to cut cost and scale your business in nearshoring a smart way to cut cost and scale your business in nearshoring a smart way to cut cost and scale your business in nearshoring a smart way to cut cost and scale your business in nearshoring a smart way to cut cost and scale your business in nearshoring a smart way to cut cost and scale your business in kpn completes m share buyback walmart is selling a dresser for only and it stackable design is so versatile bill creating new wheatfield pendleton zip code pass house walmart is selling a dresser for only and it stackable design is so versatile planner approve latest section of tracy hill kpn completes m share buyback how to send bulk email in gmail easily and quickly integrate reward into the saas customer journey luxury yacht company continues to sail smoothly despite trump trade tariff partner with palm beach hotel nfl fining more than player for selling super bowl ticket above face value visa inc nysev share sold by panoramic investment advisor llc foreign investor in charlotte eagle lake jf project obtain conditional green card approval remy cointreau invests in nonalc spirit brand jnpr elbow to the face what canada risk by embracing economic nationalism special council meeting tuesday on capital improvement priority when to apply for social security jeff bezos is considering major cable channel purchase report trump float rebate for american thanks to tariff revenue berger montague pc investigates security claim against rxsight inc nasdaq rxst trump say he considering rebate check for american based on tariff expert discus state of american manufacturing at daily caller live event apple share edge higher ahead of q earnings release from epstein to exponential debt sign of a system unraveling openais gpt expected next month could steal the vibecoding crown from anthropic are healthcare customer service job at risk in ai automation era can you close on a heloc in less than two week here what to know now education department say it will release billion in remaining withheld grant money for school education department say it will release billion in remaining withheld grant money for school trump administration approves controversial acquisition of paramount by skydance medium discovery announces detail of second quarter conference call and webcast discovery announces detail of second quarter conference call and webcast behind trump war with powell a battering ram with million follower mcdonalds to test drink from shuttering cosmcs spinoff at some restaurant longbow finance sa acquires share of abbvie inc nyseabbv prospect financial group llc acquires share of abbvie inc nyseabbv hedeker wealth llc trim stake in vanguard sp etf nysearcavoo longbow finance sa buy share of abbvie inc nyseabbv longbow finance sa boost stock holding in intuitive surgical inc nasdaqisrg curb the splurge prominent chiropractor accused of sending indecent material cbnnws plan carbondale cheapening genocide strip the nea of it charter and other commentary uk germany france call for immediate ceasefire in gaza yurachek alphabet inc nasdaqgoog share sold by elefante mark b nicholas wealth llc trim stock position in alphabet inc nasdaqgoog ice begin deportation flight from florida alligator alcatraz detention facility elon musk hit back at trump after president said he didnt want to strip tesla bos of subsidy ap news summary at pm edt sylacauga voter can hear their candidate idea aug trump official hit with fascist accusation over stephen colbert ouster pharaoh maharaja and the making of a multipolar world house ethic committee order aoc to pay additional for met gala showing countryrap star gain upper hand in custody battle over daughter trump team want to rewrite history by axing book on slavery letter to the editor wallis for black peterson nichols in longview race dozen of kid and adult in gaza have starved to death in july a hunger surge dozen of kid and adult in gaza have starved to death in july a hunger surge park county man allegedly plotted mass shooting at chico hot spring dozen of kid and adult in gaza have starved to death in july a hunger surge dozen of kid and adult in gaza have starved to death in july a hunger surge dozen of kid and adult in gaza have starved to death in july a hunger surge dozen of kid and adult in gaza have starved to death in july a hunger surge dozen of kid and adult in gaza have starved to death in july a hunger surge woodsonia ready to proceed with grand island good life application president say illinois state fair upping security measure this year dozen of kid and adult in gaza have starved to death in july a hunger surge dozen of kid and adult in gaza have starved to death in july a hunger surge dozen of kid and adult in gaza have starved to death in july a hunger surge illinois state fair upping security measure this year dozen of kid and adult in gaza have starved to death in july a hunger surge dozen of kid and adult in gaza have starved to death in july a hunger surge dozen of kid and adult in gaza have starved to death in july a hunger surge dozen of kid and adult in gaza have starved to death in july a hunger surge dozen of kid and adult in gaza have starved to death in july a hunger surge dozen of kid and adult in gaza have starved to death in july a hunger surge dozen of kid and adult in gaza have starved to death in july a hunger surge dozen of kid and adult in gaza have starved to death in july a hunger surge gordon melville mackenzie dozen of kid and adult in gaza have starved to death in july a hunger surge dozen of kid and adult in gaza have starved to death in july a hunger surge central florida daca recipient released from alligator alcatraz attorney say he still jailed trump sign executive order clamping down on homeless trump sign executive order clamping down on homeless nfl fining more than player for selling super bowl ticket above face value caterpillar inc nysecat share sold by ae wealth management llc hospital chaplain fired after supporting former coworker detained by ice hospital chaplain fired after supporting former coworker detained by ice hospital chaplain fired after supporting former coworker detained by ice hospital chaplain fired after supporting former coworker detained by ice hospital chaplain fired after supporting former coworker detained by ice house ethic panel tell aoc to pay another for met gala dress ucsf worker on strike over layoff in epstein furor trump struggle to shake off a controversy his ally once stoked in epstein furor trump struggle to shake off a controversy his ally once stoked in epstein furor trump struggle to shake off a controversy his ally once stoked in epstein furor trump struggle to shake off a controversy his ally once stoked in epstein furor trump struggle to shake off a controversy his ally once stoked in epstein furor trump struggle to shake off a controversy his ally once stoked in epstein furor trump struggle to shake off a controversy his ally once stoked in epstein furor trump struggle to shake off a controversy his ally once stoked in epstein furor trump struggle to shake off a controversy his ally once stoked in epstein furor trump struggle to shake off a controversy his ally once stoked in epstein furor trump struggle to shake off a controversy his ally once stoked in epstein furor trump struggle to shake off a controversy his ally once stoked in epstein furor trump struggle to shake off a controversy his ally once stoked in epstein furor trump struggle to shake off a controversy his ally once stoked in epstein furor trump struggle to shake off a controversy his ally once stoked in epstein furor trump struggle to shake off a controversy his ally once stoked in epstein furor trump struggle to shake off a controversy his ally once stoked in epstein furor trump struggle to shake off a controversy his ally once stoked in epstein furor trump struggle to shake off a controversy his ally once stoked in epstein furor trump struggle to shake off a controversy his ally once stoked in epstein furor trump struggle to shake off a controversy his ally once stoked in epstein furor trump struggle to shake off a controversy his ally once stoked in epstein furor trump struggle to shake off a controversy his ally once stoked in epstein furor trump struggle to shake off a controversy his ally once stoked in epstein furor trump struggle to shake off a controversy his ally once stoked in epstein furor trump struggle to shake off a controversy his ally once stoked in epstein furor trump struggle to shake off a controversy his ally once stoked in epstein furor trump struggle to shake off a controversy his ally once stoked in epstein furor trump struggle to shake off a controversy his ally once stoked in epstein furor trump struggle to shake off a controversy his ally once stoked in epstein furor trump struggle to shake off a controversy his ally once stoked in epstein furor trump struggle to shake off a controversy his ally once stoked in epstein furor trump struggle to shake off a controversy his ally once stoked in epstein furor trump struggle to shake off a controversy his ally once stoked in epstein furor trump struggle to shake off a controversy his ally once stoked in epstein furor trump struggle to shake off a controversy his ally once stoked in epstein furor trump struggle to shake off a controversy his ally once stoked in epstein furor trump struggle to shake off a controversy his ally once stoked in epstein furor trump struggle to shake off a controversy his ally once stoked in epstein furor trump struggle to shake off a controversy his ally once stoked in epstein furor trump struggle to shake off a controversy his ally once stoked in epstein furor trump struggle to shake off a controversy his ally once stoked in epstein furor trump struggle to shake off a controversy his ally once stoked in epstein furor trump struggle to shake off a controversy his ally once stoked in epstein furor trump struggle to shake off a controversy his ally once stoked in epstein furor trump struggle to shake off a controversy his ally once stoked in epstein furor trump struggle to shake off a controversy his ally once stoked in epstein furor trump struggle to shake off a controversy his ally once stoked in epstein furor trump struggle to shake off a controversy his ally once stoked in epstein furor trump struggle to shake off a controversy his ally once stoked in epstein furor trump struggle to shake off a controversy his ally once stoked in epstein furor trump struggle to shake off a controversy his ally once stoked ae wealth management llc raise holding in vanguard midcap etf nysearcavo ae wealth management llc sell share of general mill inc nysegis ae wealth management llc acquires share of the charles schwab corporation nyseschw burney co sell share of zoetis inc nysezts two longterm dividend stock to secure your financial future molina healthcare inc nysemoh share purchased by tempus wealth planning llc dimensional emerging core equity market etf nysearcadfae share acquired by avantax advisory service inc tempus wealth planning llc buy share of vanguard information technology etf nysearcavgt gilbert cook inc acquires new holding in ovintiv inc nyseovv ishares year treasury bond etf nasdaqtlt stake cut by avantax advisory service inc royce associate lp buy share of pattersonuti energy inc nasdaqpten first trust managed municipal etf nasdaqfmb share sold by tempus wealth planning llc royce associate lp trim stake in silvercrest asset management group inc nasdaqsamg tempus wealth planning llc ha stock holding in vanguard information technology etf nysearcavgt avantax advisory service inc grows position in oracle corporation nyseorcl avantax advisory service inc trim holding in ishares year treasury bond etf nasdaqtlt gilbert cook inc buy new stake in ovintiv inc nyseovv amazon is selling a slender smartwatch for and shopper love it accurate fitness tracking channel rewind traffic chaos growth safety upgrade video vault future file consumer not beefing about these price brookwood investment group llc make new investment in kinder morgan inc nysekmi brookwood investment group llc buy share of first trust short duration managed municipal etf nysearcafsmb brookwood investment group llc take position in cigna group nyseci brookwood investment group llc acquires share of ishares core sp smallcap etf nysearcaijr organon co nyseogn share sold by royce associate lp royce associate lp lower stake in powell industry inc nasdaqpowl brookwood investment group llc buy share of southern company the nyseso organon co nyseogn share sold by royce associate lp brookwood investment group llc invests in kinder morgan inc nysekmi royce associate lp ha million stock position in ishares russell value etf nysearcaiwn brookwood investment group llc acquires share of crown castle inc nysecci brookwood investment group llc acquires share of vanguard consumer discretionary etf nysearcavcr usjapan trade deal give trump control over billion it could be vapor ware fortune globalt investment llc ga cut stock holding in uber technology inc nyseuber nikko asset management america inc sell share of alphabet inc nasdaqgoog marietta investment partner llc invests in costco wholesale corporation nasdaqcost journey strategic wealth llc sell share of chipotle mexican grill inc nysecmg journey strategic wealth llc ha position in morgan stanley nysems share in gsk plc sponsored adr nysegsk purchased by lifework advisor llc freeportmcmoran inc nysefcx share sold by journey strategic wealth llc journey strategic wealth llc ha stock position in msci inc nysemsci brookwood investment group llc purchase share of philip morris international inc nysepm journey strategic wealth llc ha stake in vanguard smallcap value etf nysearcavbr mgm resort international nysemgm share acquired by journey strategic wealth llc advisor o llc raise position in morgan stanley nysems lifework advisor llc make new investment in steel dynamic inc nasdaqstld advisor o llc purchase share of broadridge financial solution inc nysebr lifework advisor llc sell share of spdr portfolio long term treasury etf nysearcasptl one wealth advisor llc grows position in emerson electric co nyseemr journey strategic wealth llc purchase share of msci inc nysemsci advisor o llc ha stock position in trane technology plc nysett advisor o llc sell share of charles river laboratory international inc nysecrl brookwood investment group llc buy share of philip morris international inc nysepm dimensional international value etf nysearcadfiv share sold by journey strategic wealth llc advanced micro device inc nasdaqamd stock position lifted by lifework advisor llc freeportmcmoran inc nysefcx share sold by journey strategic wealth llc journey strategic wealth llc ha position in ford motor company nysef brandywine oak private wealth llc decrease stake in international business machine corporation nyseibm hb wealth management llc grows stock position in johnson johnson nysejnj washington gas customer in virginia to see a surcharge this summer aurora facing significant hole in budget mayor say advisor o llc increase stock position in bristol myers squibb company nysebmy wec energy group inc nysewec share acquired by carnegie investment counsel advisor o llc decrease stock holding in teledyne technology incorporated nysetdy carnegie investment counsel cut stock position in invesco bulletshares high yield corporate bond etf nasdaqbsjq advisor o llc lower stock holding in becton dickinson and company nysebdx advisor o llc ha stock position in ishares msci usa quality factor etf batsqual advisor o llc buy share of arista network inc nyseanet archer daniel midland company nyseadm share purchased by advisor o llc intuit inc nasdaqintu share purchased by carnegie investment counsel adobe inc nasdaqadbe share bought by advisor o llc advisor o llc sell share of equifax inc nyseefx kovack advisor inc buy share of invesco aerospace defense etf nysearcappa fair isaac corporation nysefico share sold by advisor o llc advisor o llc ha stock position in ge vernova inc nysegev eversource energy nysees share bought by advisor o llc honeywell international inc nasdaqhon share acquired by advisor o llc advisor o llc acquires share of arista network inc nyseanet carnegie investment counsel invests in equifax inc nyseefx advisor o llc ha position in sei investment company nasdaqseic carnegie investment counsel ha million holding in fiserv inc nysefi advisor o llc sell share of ishares msci usa quality factor etf batsqual wedbush security inc decrease position in capital group u multisector income etf nysearcacgms phoenix financial ltd increase stock holding in barrick mining corporation nyseb carnegie investment counsel invests million in bank of montreal nysebmo machina capital sa make new investment in marvell technology inc nasdaqmrvl phoenix financial ltd invests in tmobile u inc nasdaqtmus carnegie investment counsel ha million stock position in ishares global tech etf nysearcaixn carnegie investment counsel sell share of ishares msci eafe etf nysearcaefa schwab u largecap etf nysearcaschx share acquired by carnegie investment counsel carnegie investment counsel acquires share of the walt disney company nysedis machina capital sa cut stock position in workday inc nasdaqwday carnegie investment counsel purchase share of dover corporation nysedov carnegie investment counsel sell share of comcast corporation nasdaqcmcsa phoenix financial ltd invests in tmobile u inc nasdaqtmus carnegie investment counsel acquires share of american express company nyseaxp phoenix financial ltd cut stake in synopsys inc nasdaqsnps ishares russell growth etf nysearcaiwf share bought by carnegie investment counsel carnegie investment counsel raise stock position in american electric power company inc nasdaqaep carnegie investment counsel purchase share of vanguard smallcap growth etf nysearcavbk letter bzdok pinnacle associate ltd acquires share of ishares russell midcap growth etf nysearcaiwp tempus wealth planning llc decrease stock position in first trust managed municipal etf nasdaqfmb avantax advisory service inc raise position in oracle corporation nyseorcl house democrat look to get copy of late sex offender jeffrey epstein birthday book astronomer hire gwyneth paltrow to be temporary spokesperson after kiss cam drama child who starved to death in gaza weighed less than when she wa born child who starved to death in gaza weighed less than when she wa born child who starved to death in gaza weighed less than when she wa born child who starved to death in gaza weighed less than when she wa born child who starved to death in gaza weighed less than when she wa born child who starved to death in gaza weighed less than when she wa born watch a the company that employed the coldplay cheater run an absolute clinic on damage control child who starved to death in gaza weighed less than when she wa born youth delegate represent polk county at th tennessee h congress child who starved to death in gaza weighed less than when she wa born child who starved to death in gaza weighed less than when she wa born child who starved to death in gaza weighed less than when she wa born child who starved to death in gaza weighed less than when she wa born child who starved to death in gaza weighed less than when she wa born child who starved to death in gaza weighed less than when she wa born child who starved to death in gaza weighed less than when she wa born child who starved to death in gaza weighed less than when she wa born child who starved to death in gaza weighed less than when she wa born child who starved to death in gaza weighed less than when she wa born child who starved to death in gaza weighed less than when she wa born child who starved to death in gaza weighed less than when she wa born child who starved to death in gaza weighed less than when she wa born child who starved to death in gaza weighed less than when she wa born child who starved to death in gaza weighed less than when she wa born child who starved to death in gaza weighed less than when she wa born child who starved to death in gaza weighed less than when she wa born child who starved to death in gaza weighed less than when she wa born child who starved to death in gaza weighed less than when she wa born child who starved to death in gaza weighed less than when she wa born child who starved to death in gaza weighed less than when she wa born child who starved to death in gaza weighed less than when she wa born child who starved to death in gaza weighed less than when she wa born child who starved to death in gaza weighed less than when she wa born child who starved to death in gaza weighed less than when she wa born child who starved to death in gaza weighed less than when she wa born child who starved to death in gaza weighed less than when she wa born child who starved to death in gaza weighed less than when she wa born child who starved to death in gaza weighed less than when she wa born child who starved to death in gaza weighed less than when she wa born child who starved to death in gaza weighed less than when she wa born child who starved to death in gaza weighed less than when she wa born child who starved to death in gaza weighed less than when she wa born child who starved to death in gaza weighed less than when she wa born child who starved to death in gaza weighed less than when she wa born child who starved to death in gaza weighed less than when she wa born child who starved to death in gaza weighed less than when she wa born child who starved to death in gaza weighed less than when she wa born share in the kroger co nysekr acquired by brookwood investment group llc royce associate lp ha million stock holding in ishares russell value etf nysearcaiwn brookwood investment group llc make new investment in crown castle inc nysecci clash over a promotion put hegseth at odds with his general horror upon horror how u congress responded to mass hunger in gaza lifework advisor llc boost stake in advanced micro device inc nasdaqamd journey strategic wealth llc ha stock position in ford motor company nysef drugpanel member who resigned lament ignoring of expert advice role of politics bc eaton corporation plc nyseetn share sold by advisor o llc carnegie investment counsel boost holding in cdw corporation nasdaqcdw sei investment company nasdaqseic stock position boosted by advisor o llc carnegie investment counsel ha million position in fiserv inc nysefi advisor o llc purchase share of eog resource inc nyseeog ap news summary at pm edt ap news summary at pm edt woman in legislature across the u fight for potty parity woman in legislature across the u fight for potty parity woman in legislature across the u fight for potty parity woman in legislature across the u fight for potty parity woman in legislature across the u fight for potty parity woman in legislature across the u fight for potty parity woman in legislature across the u fight for potty parity woman in legislature across the u fight for potty parity president donald trump say japan will invest billion in u at his direction it may not be a sure thing forum epstein case dont look away from this controversy letter president trump vision for america will have calamitous consequence woman in legislature across the u fight for potty parity quentin and carnell lake team up to tackle heart valve disease ishares russell growth etf nysearcaiwf share acquired by carnegie investment counsel phoenix financial ltd sell share of synopsys inc nasdaqsnps carnegie investment counsel trim stock position in booking holding inc nasdaqbkng carnegie investment counsel purchase share of american express company nyseaxp carnegie investment counsel increase stake in american electric power company inc nasdaqaep omnicom group inc nyseomc share sold by carnegie investment counsel analyst ending digital service tax strict regulation must be part of eu trade deal analyst ending digital service tax strict regulation must be part of eu trade deal analyst ending digital service tax strict regulation must be part of eu trade deal analyst ending digital service tax strict regulation must be part of eu trade deal analyst ending digital service tax strict regulation must be part of eu trade deal analyst ending digital service tax strict regulation must be part of eu trade deal analyst ending digital service tax strict regulation must be part of eu trade deal analyst ending digital service tax strict regulation must be part of eu trade deal analyst ending digital service tax strict regulation must be part of eu trade deal analyst ending digital service tax strict regulation must be part of eu trade deal analyst ending digital service tax strict regulation must be part of eu trade deal analyst ending digital service tax strict regulation must be part of eu trade deal amazon backtoschool apple sale is live airpods ipads and macbooks from your source for insightful business news and analysis bedel financial consulting inc sell share of vanguard total stock market etf nysearcavti nisource inc nyseni share sold by spectrum wealth counsel llc kapitalo investimentos ltda purchase share of flower food inc nyseflo spectrum wealth counsel llc buy share of verizon communication inc nysevz spectrum wealth counsel llc purchase new share in ishares russell value etf nysearcaiwd spectrum wealth counsel llc trim position in ishares msci eafe growth etf batsefg kapitalo investimentos ltda acquires share of hess corporation nysehes exelon corporation nasdaqexc share sold by kapitalo investimentos ltda kapitalo investimentos ltda lower position in fox corporation nasdaqfoxa kapitalo investimentos ltda raise stock position in nrg energy inc nysenrg bedel financial consulting inc increase stock holding in pacer u cash cow etf batscowz kapitalo investimentos ltda buy share of taiwan semiconductor manufacturing company ltd nysetsm spectrum wealth counsel llc ha million stake in illinois tool work inc nyseitw kapitalo investimentos ltda purchase share of regency center corporation nasdaqreg trump announces trade agreement with eu to set tariff at ask eartha how can home energy assessment help save money charge because i unplugged a cord paris la vega hotel tiny print fee sum up why visitor are fleeing the city view from the wing nextera energy inc nysenee share bought by mrp capital investment llc advantage trust co decrease stake in first trust smid cap rising dividend achiever etf nasdaqsdvy customer data stolen at u branch of german insurer allianz carmax inc nysekmx share bought by assetmark inc ap news summary at pm edt ap news summary at pm edt ap news summary at pm edt urgent recall of popular cheeseburger that can trigger lifethreatening allergic reaction daily mail bay bank of virginia otcmktsbayk and renasant nasdaqrnst critical contrast intellus advisor llc reduces stock position in raymond james financial inc nyserjf intellus advisor llc ha million holding in alphabet inc nasdaqgoog applied finance capital management llc grows stock position in cardinal health inc nysecah otis worldwide corporation nyseotis share acquired by empirical finance llc intellus advisor llc acquires share of tesla inc nasdaqtsla paccar inc nasdaqpcar receives average recommendation of moderate buy from brokerage ishares silver trust nysearcaslv share bought by intellus advisor llc avantax advisory service inc buy share of servicenow inc nysenow what to know about soda sweetener a sugar return to american coke what to know about soda sweetener a sugar return to american coke what to know about soda sweetener a sugar return to american coke what to know about soda sweetener a sugar return to american coke what to know about soda sweetener a sugar return to american coke what to know about soda sweetener a sugar return to american coke what to know about soda sweetener a sugar return to american coke what to know about soda sweetener a sugar return to american coke what to know about soda sweetener a sugar return to american coke what to know about soda sweetener a sugar return to american coke what to know about soda sweetener a sugar return to american coke walmart is selling a weatherresistant deck box for and shopper say it storage space is impressive what to know about soda sweetener a sugar return to american coke what to know about soda sweetener a sugar return to american coke what to know about soda sweetener a sugar return to american coke what to know about soda sweetener a sugar return to american coke what to know about soda sweetener a sugar return to american coke what to know about soda sweetener a sugar return to american coke what to know about soda sweetener a sugar return to american coke what to know about soda sweetener a sugar return to american coke what to know about soda sweetener a sugar return to american coke what to know about soda sweetener a sugar return to american coke what to know about soda sweetener a sugar return to american coke what to know about soda sweetener a sugar return to american coke what to know about soda sweetener a sugar return to american coke accretive wealth partner llc take position in shell plc unsponsored adr nyseshel ark investment management llc ha million holding in ptc inc nasdaqptc delta air line inc nysedal share acquired by empirical finance llc ark investment management llc sell share of faro technology inc nasdaqfaro ark investment management llc ha million holding in'
